# Adds a new "Statement of truth" heading paragraph and a new content
# table immediately after the final table in the document (the
# "Defendant's legal representative" details table that ends with the
# "DX address <<respondent.representative.dxAddress>>" cell), and before
# the trailing empty paragraph / section break.

$d = $word.ActiveDocument

# Locate the last table in the document body (the one whose last cell
# contains the DX address merge field) and get a collapsed Range right
# after it.
$tbl = $d.Tables.Item($d.Tables.Count)
$insertRange = $d.Range($tbl.Range.End, $tbl.Range.End)

$fragmentXml = @'
<w:p>
  <w:pPr>
    <w:pBdr>
      <w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/>
    </w:pBdr>
    <w:spacing w:before="120" w:after="120"/>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:lang w:eastAsia="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="GDSTransportWebsite" w:hAnsi="GDSTransportWebsite"/>
      <w:color w:val="0A0A0A"/>
      <w:sz w:val="36"/>
      <w:szCs w:val="36"/>
    </w:rPr>
    <w:t>Statement of truth</w:t>
  </w:r>
</w:p>
<w:tbl>
  <w:tblPr>
    <w:tblStyle w:val="TableGrid"/>
    <w:tblW w:w="10060" w:type="dxa"/>
    <w:tblBorders>
      <w:top w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:left w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:bottom w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:insideH w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:insideV w:val="none" w:sz="0" w:space="0" w:color="auto"/>
    </w:tblBorders>
    <w:tblLayout w:type="fixed"/>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="10060"/>
  </w:tblGrid>
  <w:tr>
    <w:trPr>
      <w:trHeight w:val="388"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="10060" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:pStyle w:val="NormalWeb"/>
          <w:spacing w:before="120" w:beforeAutospacing="0" w:after="120" w:afterAutospacing="0"/>
          <w:rPr>
            <w:rFonts w:ascii="GDSTransportWebsite" w:hAnsi="GDSTransportWebsite"/>
            <w:color w:val="0A0A0A"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:ascii="GDSTransportWebsite" w:hAnsi="GDSTransportWebsite"/>
            <w:color w:val="0A0A0A"/>
          </w:rPr>
          <w:t>The claimant believes that the facts in this claim form are true.</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:rFonts w:ascii="GDSTransportWebsite" w:hAnsi="GDSTransportWebsite"/>
            <w:color w:val="0A0A0A"/>
          </w:rPr>
        </w:pPr>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:rFonts w:ascii="GDSTransportWebsite" w:hAnsi="GDSTransportWebsite"/>
            <w:color w:val="0A0A0A"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:ascii="GDSTransportWebsite" w:hAnsi="GDSTransportWebsite"/>
            <w:color w:val="0A0A0A"/>
          </w:rPr>
          <w:t>I am duly authorised by the claimant to sign this statement.</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:rFonts w:ascii="GDSTransportWebsite" w:hAnsi="GDSTransportWebsite"/>
            <w:color w:val="0A0A0A"/>
          </w:rPr>
        </w:pPr>
      </w:p>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:rFonts w:ascii="GDSTransportWebsite" w:hAnsi="GDSTransportWebsite"/>
            <w:color w:val="0A0A0A"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:ascii="GDSTransportWebsite" w:hAnsi="GDSTransportWebsite"/>
            <w:color w:val="0A0A0A"/>
          </w:rPr>
          <w:t>The claimant understands that proceedings for contempt of court may be brought against anyone who makes, or causes to be made, a false statement in a document verified by a statement of truth without an honest belief in its truth</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:ascii="GDSTransportWebsite" w:hAnsi="GDSTransportWebsite"/>
            <w:color w:val="0A0A0A"/>
          </w:rPr>
          <w:t>.</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
'@

$packageXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
$fragmentXml
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$insertRange.InsertXML($packageXml) | Out-Null
